# chaitanya basic scenarios add
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F3").Value = 899.1
# Match the number format / style of the General-formatted cells in the row
# (copy formatting only, so the value just set above is preserved).
$wsSummary.Range("C3").Copy()
$wsSummary.Range("F3").PasteSpecial(-4122)

# --- Repayment Schedule sheet ---
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
# Add an (empty) formatted cell at O2, matching its neighbour N2.
$wsRepay.Range("N2").Copy()
$wsRepay.Range("O2").PasteSpecial(-4122)

# --- Transactions sheet ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 222
$wsTrans.Range("A3").Value = 221

# --- Restore/update each sheet's remembered selection ---
$wsSummary.Range("B3").Select()
$wsRepay.Range("F7").Select()
$wsTrans.Range("B2").Select()
